$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

# "2025-11-28" and "251128" look numeric/date-like to Excel's auto-detect,
# so force Text formatting just long enough to land them as literal strings
# (matching the rest of the column), then drop the format again so the new
# row doesn't pick up a distinct cell style from the rest of the sheet.
$ws.Range("A" + $row + ":" + "E" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-11-28"
$ws.Range("B" + $row).Value = "Pick 4"
$ws.Range("C" + $row).Value = "251128"
$ws.Range("D" + $row).Value = "2-3-8-2"
$ws.Range("E" + $row).Value = "2025-11-28T21:38:20.540+04:00"

$ws.Range("A" + $row + ":" + "E" + $row).ClearFormats()
